# Auto-generated PowerShell Excel COM-interop script
# Updates cryptocurrency price/volume data (and two row name/link/price/volume swaps)
# to match the target commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to store the given string as literal text, even when it
    # looks like a number (e.g. "1.000" or "0.9992"), mirroring how the
    # original workbook stores these as inline/shared strings rather than
    # numeric values. We temporarily apply a text number format while setting
    # the value, then restore General format and clear the style so the cell
    # ends up with no special style applied (matching the source file).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '29.869.26'
Set-TextValue $ws.Range("E2") '  -0.05%  '
Set-TextValue $ws.Range("D3") '1.887.87'
Set-TextValue $ws.Range("E3") '  -0.35%  '
Set-TextValue $ws.Range("E4") '  +0.08%  '
Set-TextValue $ws.Range("D5") '0.7709'
Set-TextValue $ws.Range("E5") '  -0.63%  '
Set-TextValue $ws.Range("D6") '242.71'
Set-TextValue $ws.Range("D7") '1.000'
Set-TextValue $ws.Range("E7") '  +0.05%  '
Set-TextValue $ws.Range("D8") '0.3121'
Set-TextValue $ws.Range("E8") '  -0.84%  '
Set-TextValue $ws.Range("D9") '25.60'
Set-TextValue $ws.Range("E9") '  +0.30%  '
Set-TextValue $ws.Range("D10") '0.07168'
Set-TextValue $ws.Range("E10") '  -5.20%  '
Set-TextValue $ws.Range("D11") '0.08611'
Set-TextValue $ws.Range("E11") '  +6.24%  '
Set-TextValue $ws.Range("D12") '0.7637'
Set-TextValue $ws.Range("E12") '  -1.24%  '
Set-TextValue $ws.Range("D13") '1.903.38'
Set-TextValue $ws.Range("E13") '  -2.01%  '
Set-TextValue $ws.Range("D14") '5.359'
Set-TextValue $ws.Range("E14") '  -2.82%  '
Set-TextValue $ws.Range("D15") '93.53'
Set-TextValue $ws.Range("E15") '  +1.14%  '
Set-TextValue $ws.Range("D16") '6.155'
Set-TextValue $ws.Range("E16") '  -1.85%  '
Set-TextValue $ws.Range("D17") '29.812.52'
Set-TextValue $ws.Range("E17") '  -0.10%  '
Set-TextValue $ws.Range("D18") '13.75'
Set-TextValue $ws.Range("E18") '  -1.87%  '
Set-TextValue $ws.Range("D19") '244.25'
Set-TextValue $ws.Range("E19") '  -0.20%  '
Set-TextValue $ws.Range("D20") '0.000007799'
Set-TextValue $ws.Range("E20") '  -1.36%  '
Set-TextValue $ws.Range("D21") '2.152.86'
Set-TextValue $ws.Range("E21") '  +1.06%  '
Set-TextValue $ws.Range("D22") '0.9992'
Set-TextValue $ws.Range("E22") '  -0.09%  '
Set-TextValue $ws.Range("D23") '8.009'
Set-TextValue $ws.Range("E23") '  -1.86%  '
Set-TextValue $ws.Range("E24") '  +0.12%  '
Set-TextValue $ws.Range("D25") '0.1652'
Set-TextValue $ws.Range("E25") '  +5.13%  '
Set-TextValue $ws.Range("D26") '9.366'
Set-TextValue $ws.Range("E26") '  -1.11%  '
Set-TextValue $ws.Range("D27") '162.45'
Set-TextValue $ws.Range("E27") '  -0.19%  '
Set-TextValue $ws.Range("D28") '18.73'
Set-TextValue $ws.Range("E28") '  -0.30%  '
Set-TextValue $ws.Range("D29") '2.033'
Set-TextValue $ws.Range("E29") '  -0.66%  '
Set-TextValue $ws.Range("D30") '1.459'
Set-TextValue $ws.Range("E30") '  +1.46%  '
Set-TextValue $ws.Range("D31") '1.534'
Set-TextValue $ws.Range("E31") '  -1.22%  '
Set-TextValue $ws.Range("D32") '4.504'
Set-TextValue $ws.Range("E32") '  +0.48%  '
Set-TextValue $ws.Range("D33") '4.095'
Set-TextValue $ws.Range("E33") '  -0.14%  '
Set-TextValue $ws.Range("D34") '0.05435'
Set-TextValue $ws.Range("E34") '  -1.72%  '
Set-TextValue $ws.Range("D35") '1.240'
Set-TextValue $ws.Range("E35") '  -1.60%  '
Set-TextValue $ws.Range("D36") '0.7424'
Set-TextValue $ws.Range("E36") '  -2.01%  '
Set-TextValue $ws.Range("D37") '0.9989'
Set-TextValue $ws.Range("E37") '  -0.03%  '
Set-TextValue $ws.Range("D38") '2.700'
Set-TextValue $ws.Range("E38") '  +2.14%  '
Set-TextValue $ws.Range("D39") '0.01955'
Set-TextValue $ws.Range("E39") '  +1.25%  '
Set-TextValue $ws.Range("D40") '2.781'
Set-TextValue $ws.Range("E40") '  -0.30%  '
Set-TextValue $ws.Range("E41") '  +0.48%  '
Set-TextValue $ws.Range("D42") '1.108.15'
Set-TextValue $ws.Range("E42") '  -4.16%  '
Set-TextValue $ws.Range("D43") '72.97'
Set-TextValue $ws.Range("E43") '  -1.49%  '
Set-TextValue $ws.Range("D44") '6.065'
Set-TextValue $ws.Range("E44") '  +2.12%  '
Set-TextValue $ws.Range("D45") '0.8510'
Set-TextValue $ws.Range("E45") '  +0.30%  '
Set-TextValue $ws.Range("D46") '1.000'
Set-TextValue $ws.Range("E46") '  +0.01%  '
Set-TextValue $ws.Range("D47") '102.22'
Set-TextValue $ws.Range("E47") '  -0.25%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D48") '1.863'
Set-TextValue $ws.Range("E48") '  -1.94%  '
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D49") '7.640'
Set-TextValue $ws.Range("E49") '  +1.68%  '
$ws.Range("B50").Value = 'SynthetixNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
Set-TextValue $ws.Range("D50") '2.998'
Set-TextValue $ws.Range("E50") '  -4.06%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range("D51") '2.047.69'
Set-TextValue $ws.Range("E51") '  +0.52%  '
